$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -6
$ws.Range("F7").Value = -5
$ws.Range("F8").Value = -3
$ws.Range("F9").Value = -6
$ws.Range("F10").Value = -11
$ws.Range("F16").Value = -3
$ws.Range("F18").Value = -5
$ws.Range("F19").Value = -2
$ws.Range("F21").Value = 0
$ws.Range("F24").Value = 1
